$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-TEMP_SENS_1")

# Insert a new blank row above the current row 6 ("Mounting hole, 4.2mm"),
# pushing that and all subsequent rows down by one.
$ws.Rows.Item(6).Insert()

# Row 5: connector header changed from right-angle 1.25mm to vertical 2.54mm.
$ws.Cells.Item(5, 1).Value = "CONN HEADER VERT 3POS 2.54MM"

# New row 6: single layer pad test points JP2 & JP3.
$ws.Cells.Item(6, 1).Value = "Single layer pad TP"
$ws.Cells.Item(6, 2).Value = "JP2, JP3"
$ws.Cells.Item(6, 3).Value = 2

$wb.Save()
